$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns remain text, matching the
# source data which stores these as formatted strings (not numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "69.373.65"
$ws.Range("E2").Value = "  -2.41%  "

$ws.Range("D3").Value = "3.537.95"
$ws.Range("E3").Value = "  -4.09%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "582.15"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").Value = "171.91"
$ws.Range("E6").Value = "  -3.48%  "

$ws.Range("D7").Value = "3.531.60"
$ws.Range("E7").Value = "  -3.92%  "

$ws.Range("E8").Value = "  -1.01%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("E10").Value = "  -4.85%  "

$ws.Range("D11").Value = "6.69"
$ws.Range("E11").Value = "  -2.39%  "

$ws.Range("D12").Value = "0.587"
$ws.Range("E12").Value = "  -4.01%  "

$ws.Range("D13").Value = "47.61"
$ws.Range("E13").Value = "  -3.14%  "

$ws.Range("E14").Value = "  -4.47%  "

$ws.Range("D15").Value = "4.096.44"
$ws.Range("E15").Value = "  -4.32%  "

$ws.Range("E16").Value = "  -5.27%  "

$ws.Range("D17").Value = "627.45"
$ws.Range("E17").Value = "  -7.63%  "

$ws.Range("D18").Value = "3.538.47"
$ws.Range("E18").Value = "  -4.07%  "

$ws.Range("D19").Value = "69.351.28"
$ws.Range("E19").Value = "  -2.56%  "

$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("D21").Value = "17.54"
$ws.Range("E21").Value = "  -2.45%  "

$ws.Range("D22").Value = "11.23"
$ws.Range("E22").Value = "  -3.11%  "

$ws.Range("E23").Value = "  -5.32%  "

$ws.Range("D24").Value = "16.06"
$ws.Range("E24").Value = "  -7.60%  "

$ws.Range("D25").Value = "97.90"
$ws.Range("E25").Value = "  -3.96%  "

$ws.Range("E26").Value = "  -3.93%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("E28").Value = "  -6.85%  "

$ws.Range("D29").Value = "9.32"
$ws.Range("E29").Value = "  -8.99%  "

$ws.Range("D30").Value = "32.95"
$ws.Range("E30").Value = "  -6.31%  "

$ws.Range("E31").Value = "  -7.63%  "

$ws.Range("E32").Value = "  -6.32%  "

$ws.Range("E33").Value = "  -6.16%  "

$ws.Range("D34").Value = "7.01"
$ws.Range("E34").Value = "  -7.53%  "

$ws.Range("D35").Value = "637.38"
$ws.Range("E35").Value = "  +9.75%  "

$ws.Range("E36").Value = "  -3.69%  "

$ws.Range("D37").Value = "3.51"
$ws.Range("E37").Value = "  -13.29%  "

$ws.Range("E38").Value = "  -4.66%  "

$ws.Range("D39").Value = "57.47"
$ws.Range("E39").Value = "  -2.18%  "

$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("E41").Value = "  -1.38%  "

$ws.Range("E42").Value = "  -5.19%  "

$ws.Range("D43").Value = "3.395.04"
$ws.Range("E43").Value = "  -6.39%  "

$ws.Range("E44").Value = "  -5.98%  "

$ws.Range("E45").Value = "  -6.64%  "

$ws.Range("D46").Value = "0.0₃0703"
$ws.Range("E46").Value = "  -8.51%  "

$ws.Range("D47").Value = "2.57"
$ws.Range("E47").Value = "  -6.99%  "

$ws.Range("E48").Value = "  -4.17%  "

$ws.Range("E49").Value = "  -2.40%  "

$ws.Range("D50").Value = "5.74"

$ws.Range("D51").Value = "132.15"
$ws.Range("E51").Value = "  -1.74%  "
